$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("Created _", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph text 'Created _'"
}
$r.Collapse(0)

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7A7E0192" w14:textId="6BA83071" w:rsidR="00353029" w:rsidRDefault="00353029" w:rsidP="006D2233"><w:pPr><w:pStyle w:val="Textkrper"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Created _</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>mobile.scss</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> without menu branding &amp; nav</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Edited the menu screen, had some problems with the picture</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>31.12.2022</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Fixed the picture problem partly, changed the text and link colors</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Googled how to make the picture in the menu responsive (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>didngt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> work)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I redesigned the About Me page and added some additional background images</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
Write-Output "done"
